# Map073.xlsx "huge v2 update" patch:
#   - Column B becomes a full translation column mirroring column A for
#     every row (A1:A231) that doesn't already carry a value in B.
#   - The stray D143 cell (a leftover duplicate of the A/B143 string) is
#     removed, which collapses the sheet's used range back down from
#     A1:D231 to A1:B231.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 231

for ($r = 1; $r -le $lastRow; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $bVal = $bCell.Value()
    if ([string]::IsNullOrEmpty($bVal)) {
        $aCell = $ws.Cells.Item($r, 1)
        $aCell.Copy($bCell)
    }
}

# Drop the extra column D cell that only existed on row 143; this also
# shrinks the sheet dimension from A1:D231 down to A1:B231.
$ws.Range("D143").ClearContents()
